$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 288.4375
$ws.Range("I2").Value = 151.66667
$ws.Range("J2").Value = 464.2857
$ws.Range("K2").Value = 151.66667
$ws.Range("L2").Value = 464.2857
$ws.Range("M2").Value = -38.66667000000001
$ws.Range("N2").Value = -690.2857
$ws.Range("H4").Value = 1215.6666
$ws.Range("I4").Value = 1167.625
$ws.Range("J4").Value = 1600
$ws.Range("K4").Value = 1167.625
$ws.Range("L4").Value = 1600
$ws.Range("M4").Value = -1053.625
$ws.Range("N4").Value = -1828
$ws.Range("H5").Value = 337.9
$ws.Range("I5").Value = 213.16667
$ws.Range("J5").Value = 525
$ws.Range("K5").Value = 213.16667
$ws.Range("L5").Value = 525
$ws.Range("M5").Value = -98.16667000000001
$ws.Range("N5").Value = -755
$ws.Range("H28").Value = 449.66666
$ws.Range("I28").Value = 465
$ws.Range("J28").Value = 407.5
$ws.Range("K28").Value = 465
$ws.Range("L28").Value = 407.5
$ws.Range("M28").Value = 20
$ws.Range("N28").Value = -1377.5
$ws.Range("H40").Value = 2771.25
$ws.Range("I40").Value = 1885.6
$ws.Range("J40").Value = 3066.4666
$ws.Range("K40").Value = 1885.6
$ws.Range("L40").Value = 3066.4666
$ws.Range("M40").Value = -1710.6
$ws.Range("N40").Value = -3416.4666
$ws.Range("H43").Value = 1035.4706
$ws.Range("J43").Value = 700.25
$ws.Range("L43").Value = 700.25
$ws.Range("N43").Value = -838.25
$ws.Range("H51").Value = 10106673
$ws.Range("I51").Value = 22733022
$ws.Range("J51").Value = 5594.4
$ws.Range("K51").Value = 22733022
$ws.Range("L51").Value = 5594.4
$ws.Range("M51").Value = -22732538
$ws.Range("N51").Value = -6562.4
$ws.Range("H58").Value = 1566.1538
$ws.Range("I58").Value = 68.75
$ws.Range("J58").Value = 3962
$ws.Range("K58").Value = 206.25
$ws.Range("L58").Value = 11886
$ws.Range("M58").Value = -56.25
$ws.Range("N58").Value = -12186
$ws.Range("H62").Value = 19492.572
$ws.Range("I62").Value = 5449.533
$ws.Range("J62").Value = 54600.168
$ws.Range("K62").Value = 5449.533
$ws.Range("L62").Value = 54600.168
$ws.Range("M62").Value = -4825.533
$ws.Range("N62").Value = -55848.168
$ws.Range("H65").Value = 19492.572
$ws.Range("I65").Value = 5449.533
$ws.Range("J65").Value = 54600.168
$ws.Range("K65").Value = 27247.665
$ws.Range("L65").Value = 273000.84
$ws.Range("M65").Value = -24127.665
$ws.Range("N65").Value = -279240.84
$ws.Range("H74").Value = 4421
$ws.Range("I74").Value = 2303
$ws.Range("J74").Value = 5480
$ws.Range("K74").Value = 2303
$ws.Range("L74").Value = 5480
$ws.Range("M74").Value = -1367
$ws.Range("N74").Value = -7352
$ws.Range("H77").Value = 4421
$ws.Range("I77").Value = 2303
$ws.Range("J77").Value = 5480
$ws.Range("K77").Value = 11515
$ws.Range("L77").Value = 27400
$ws.Range("M77").Value = -6835
$ws.Range("N77").Value = -36760
$ws.Range("H100").Value = 2954.2856
$ws.Range("I100").Value = 2526.6667
$ws.Range("J100").Value = 3275
$ws.Range("K100").Value = 2526.6667
$ws.Range("L100").Value = 3275
$ws.Range("M100").Value = -1985.6667
$ws.Range("N100").Value = -4357
$ws.Range("H106").Value = 2285.7058
$ws.Range("I106").Value = 3400.2
$ws.Range("J106").Value = 2093.5518
$ws.Range("K106").Value = 3400.2
$ws.Range("L106").Value = 2093.5518
$ws.Range("M106").Value = -2769.2
$ws.Range("N106").Value = -3355.5518
$ws.Range("H107").Value = 474.54166
$ws.Range("I107").Value = 483.41177
$ws.Range("J107").Value = 453
$ws.Range("K107").Value = 483.41177
$ws.Range("L107").Value = 453
$ws.Range("M107").Value = 1436.58823
$ws.Range("N107").Value = -4293
$ws.Range("H113").Value = 3507.25
$ws.Range("I113").Value = 2475
$ws.Range("J113").Value = 3765.3125
$ws.Range("K113").Value = 2475
$ws.Range("L113").Value = 3765.3125
$ws.Range("M113").Value = 779
$ws.Range("N113").Value = -10273.3125
$ws.Range("H116").Value = 4533.1304
$ws.Range("I116").Value = 2897.6
$ws.Range("J116").Value = 5791.231
$ws.Range("K116").Value = 2897.6
$ws.Range("L116").Value = 5791.231
$ws.Range("M116").Value = 544.4000000000001
$ws.Range("N116").Value = -12675.231
$ws.Range("H125").Value = 3316.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3316.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 29848.5
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -34768.5
$ws.Range("H129").Value = 1401.3214
$ws.Range("J129").Value = 1678.6364
$ws.Range("L129").Value = 5035.9092
$ws.Range("N129").Value = -15035.9092
$ws.Range("H137").Value = 1898.0968
$ws.Range("J137").Value = 1512.5
$ws.Range("L137").Value = 4537.5
$ws.Range("N137").Value = -9637.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3292.1333
$ws.Range("I45").Value = 1869.2222
$ws.Range("K45").Value = 1869.2222
$ws.Range("M45").Value = -1492.2222
$ws.Range("H61").Value = 2416.6487
$ws.Range("I61").Value = 1385.8182
$ws.Range("J61").Value = 3928.5334
$ws.Range("K61").Value = 1385.8182
$ws.Range("L61").Value = 3928.5334
$ws.Range("M61").Value = -1173.8182
$ws.Range("N61").Value = -4352.5334
$ws.Range("H63").Value = 2696.9697
$ws.Range("I63").Value = 2007.6923
$ws.Range("K63").Value = 2007.6923
$ws.Range("M63").Value = -1321.6923
$ws.Range("H66").Value = 2696.9697
$ws.Range("I66").Value = 2007.6923
$ws.Range("K66").Value = 10038.4615
$ws.Range("M66").Value = -6606.461499999999
$ws.Range("H74").Value = 1621.9459
$ws.Range("I74").Value = 1530.48
$ws.Range("K74").Value = 1530.48
$ws.Range("M74").Value = -656.48
$ws.Range("H77").Value = 1621.9459
$ws.Range("I77").Value = 1530.48
$ws.Range("K77").Value = 7652.4
$ws.Range("M77").Value = -3284.4
$ws.Range("H97").Value = 718.913
$ws.Range("I97").Value = 527.25
$ws.Range("J97").Value = 1996.6666
$ws.Range("K97").Value = 527.25
$ws.Range("L97").Value = 1996.6666
$ws.Range("M97").Value = -31.25
$ws.Range("N97").Value = -2988.6666
$ws.Range("H102").Value = 2019.3125
$ws.Range("I102").Value = 1900.8182
$ws.Range("J102").Value = 2280
$ws.Range("K102").Value = 1900.8182
$ws.Range("L102").Value = 2280
$ws.Range("M102").Value = -278.8181999999999
$ws.Range("N102").Value = -5524
$ws.Range("H122").Value = 7594
$ws.Range("I122").Value = 1992.5
$ws.Range("K122").Value = 5977.5
$ws.Range("M122").Value = -3527.5
$ws.Range("H136").Value = 2416.6487
$ws.Range("I136").Value = 1385.8182
$ws.Range("J136").Value = 3928.5334
$ws.Range("K136").Value = 4157.4546
$ws.Range("L136").Value = 11785.6002
$ws.Range("M136").Value = -1607.4546
$ws.Range("N136").Value = -16885.6002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2225
$ws.Range("I99").Value = 1372.2222
$ws.Range("J99").Value = 3760
$ws.Range("K99").Value = 1372.2222
$ws.Range("L99").Value = 3760
$ws.Range("M99").Value = 125.7778000000001
$ws.Range("N99").Value = -6756
$ws.Range("H134").Value = 5982.7026
$ws.Range("I134").Value = 2885.4119
$ws.Range("J134").Value = 8615.4
$ws.Range("K134").Value = 8656.235700000001
$ws.Range("L134").Value = 25846.2
$ws.Range("M134").Value = -6121.235700000001
$ws.Range("N134").Value = -30916.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3820.25
$ws.Range("I7").Value = 6690.3335
$ws.Range("K7").Value = 6690.3335
$ws.Range("M7").Value = -6577.3335
$ws.Range("H107").Value = 577.1818
$ws.Range("I107").Value = 198.13043
$ws.Range("J107").Value = 1449
$ws.Range("K107").Value = 198.13043
$ws.Range("L107").Value = 1449
$ws.Range("M107").Value = 1721.86957
$ws.Range("N107").Value = -5289
$ws.Range("H132").Value = 3355.0833
$ws.Range("I132").Value = 2389.5557
$ws.Range("J132").Value = 3934.4
$ws.Range("K132").Value = 7168.6671
$ws.Range("L132").Value = 11803.2
$ws.Range("M132").Value = -4638.6671
$ws.Range("N132").Value = -16863.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 315.1579
$ws.Range("I5").Value = 315.1579
$ws.Range("K5").Value = 945.4737
$ws.Range("M5").Value = -833.4737
$ws.Range("H113").Value = 578.125
$ws.Range("J113").Value = 505
$ws.Range("L113").Value = 1515
$ws.Range("N113").Value = -5855
$ws.Range("H135").Value = 315.1579
$ws.Range("I135").Value = 315.1579
$ws.Range("K135").Value = 2836.4211
$ws.Range("M135").Value = -301.4211

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4890.6
$ws.Range("I70").Value = 4786.857
$ws.Range("K70").Value = 4786.857
$ws.Range("M70").Value = -4516.857
$ws.Range("H73").Value = 4890.6
$ws.Range("I73").Value = 4786.857
$ws.Range("K73").Value = 4786.857
$ws.Range("M73").Value = -3850.857
$ws.Range("H97").Value = 3520.08
$ws.Range("I97").Value = 3433.7144
$ws.Range("J97").Value = 3630
$ws.Range("K97").Value = 3433.7144
$ws.Range("L97").Value = 3630
$ws.Range("M97").Value = -2937.7144
$ws.Range("N97").Value = -4622
$ws.Range("H102").Value = 3762810
$ws.Range("H122").Value = 669437.4
$ws.Range("I122").Value = 911464.6
$ws.Range("J122").Value = 3862.5
$ws.Range("K122").Value = 2734393.8
$ws.Range("L122").Value = 11587.5
$ws.Range("M122").Value = -2731943.8
$ws.Range("N122").Value = -16487.5
$ws.Range("H132").Value = 1896403.9
$ws.Range("I132").Value = 4631832
$ws.Range("J132").Value = 2646
$ws.Range("K132").Value = 13895496
$ws.Range("L132").Value = 7938
$ws.Range("M132").Value = -13892966
$ws.Range("N132").Value = -12998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 734.3333
$ws.Range("I22").Value = 302.1
$ws.Range("J22").Value = 1598.8
$ws.Range("K22").Value = 302.1
$ws.Range("L22").Value = 1598.8
$ws.Range("M22").Value = -7.100000000000023
$ws.Range("N22").Value = -2188.8
$ws.Range("H27").Value = 734.3333
$ws.Range("I27").Value = 302.1
$ws.Range("J27").Value = 1598.8
$ws.Range("K27").Value = 302.1
$ws.Range("L27").Value = 1598.8
$ws.Range("M27").Value = -195.1
$ws.Range("N27").Value = -1812.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6758480.5
$ws.Range("I136").Value = 10870695
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 32612085
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -32609535
$ws.Range("N136").Value = -13200
